$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (G1) to the new header
# cell H1, then set its text. This mirrors the header cells B1:G1 which
# all use the bold/centered/bordered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for the "Save" column in row 2.
$ws.Range("H2").Value = 0
